$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header O2: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to columns I, J, K for data rows 3 through 23
for ($row = 3; $row -le 23; $row++) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $cell.Value = "$($cell.Text) msec"
    }
}
